{"js": "// The document ends with an empty \"Ghost\" sub-bullet list (ListParagraph,\n// numId=3) whose very last item (ilvl=0) has no text yet - it's the\n// \"Blue Cultist\" enemy entry prepared as an empty placeholder. We give that\n// paragraph its text, then append two more sub-bullets (ilvl=1) describing\n// the Blue Cultist's shot pattern and movement.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The target is the last paragraph in the document body (currently empty).\nconst items = paragraphs.items;\nconst bluecultist = items[items.length - 1];\n\n// Fill in the heading text for this enemy, then restore the run's font\n// size (inserting text resets direct character formatting) to match the\n// surrounding 18pt (sz/szCs 36 half-points) body text.\nbluecultist.insertText(\"Blue Cultist\", \"Replace\");\nawait context.sync();\nbluecultist.font.size = 18;\nbluecultist.font.sizeBidirectional = 18;\nawait context.sync();\n\n// Add the \"Shot pattern\" sub-bullet right after it, one list level deeper.\nconst shotPattern = bluecultist.insertParagraph(\n  \"Shot pattern: None, instead periodically opens two portals and reaches \" +\n  \"his two corrupted hands into them. One portal connects to directly \" +\n  \"below the player and the other somewhere nearby player. Hands deal \" +\n  \"contact damage if the player does not move. (I stole this attack from \" +\n  \"nightmare in OSRS)\",\n  \"After\"\n);\nawait context.sync();\nshotPattern.listItemOrNullObject.level = 1;\nawait context.sync();\n\n// Add the \"Movement\" sub-bullet after that, also one list level deeper.\nconst movement = shotPattern.insertParagraph(\"Movement: See red cultist\", \"After\");\nawait context.sync();\nmovement.listItemOrNullObject.level = 1;\nawait context.sync();\n", "ps1": "# The document ends with an empty \"Ghost\" sub-bullet list (ListParagraph,\n# numId=3) whose very last item (ilvl=0) has no text yet - it's the\n# \"Blue Cultist\" enemy entry prepared as an empty placeholder. We give that\n# paragraph its text, then append two more sub-bullets (ilvl=1) describing\n# the Blue Cultist's shot pattern and movement.\n\n$d = $word.ActiveDocument\n\n# The target is the last paragraph in the document body (currently empty).\n$blueCultist = $d.Paragraphs.Item($d.Paragraphs.Count)\n\n# Fill in the heading text for this enemy, then restore the run's font\n# size (setting Range.Text resets direct character formatting) to match\n# the surrounding 18pt (sz/szCs 36 half-points) body text.\n$blueCultist.Range.Text = \"Blue Cultist\"\n$blueCultist.Range.Font.Size = 18\n$blueCultist.Range.Font.SizeBi = 18\n\n# Add the \"Shot pattern\" sub-bullet right after it, one list level deeper\n# (ListLevelNumber is 1-based, so level 2 == w:ilvl val=\"1\").\n$blueCultist.Range.InsertParagraphAfter()\n$shotPattern = $d.Paragraphs.Item($d.Paragraphs.Count)\n$shotPattern.Range.Text = \"Shot pattern: None, instead periodically opens two portals and reaches his two corrupted hands into them. One portal connects to directly below the player and the other somewhere nearby player. Hands deal contact damage if the player does not move. (I stole this attack from nightmare in OSRS)\"\n$shotPattern.Range.ListFormat.ListLevelNumber = 2\n$shotPattern.Range.Font.Size = 18\n$shotPattern.Range.Font.SizeBi = 18\n\n# Add the \"Movement\" sub-bullet after that, also one list level deeper.\n$shotPattern.Range.InsertParagraphAfter()\n$movement = $d.Paragraphs.Item($d.Paragraphs.Count)\n$movement.Range.Text = \"Movement: See red cultist\"\n$movement.Range.ListFormat.ListLevelNumber = 2\n$movement.Range.Font.Size = 18\n$movement.Range.Font.SizeBi = 18\n"}
